# Applies the "FollowUpId" column changes to the GNH Card Data workbook.
#
# Summary of the change (from commit "[FIXED] several issues with parsing,
# array sizes, etc [CHANGED] imported card data to work with"):
#   - The header in H1 is renamed from "Follow Up ID" to "FollowUpId".
#   - Every "regular" card row (everything except the three FollowUp cards
#     in rows 27-29, and the rows that already reference a follow-up card
#     in rows 4, 7 and 28) gets an explicit FollowUpId value of 0 in
#     column H, marking that the card has no follow-up card attached.
#   - The selected cell / scroll position of the sheet changes as a side
#     effect of the editing session (H28 ends up selected, and the view
#     is scrolled back to the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cards")

# Rename the "Follow Up ID" header to "FollowUpId".
$ws.Range("H1").Value = "FollowUpId"

# Rows that need a default FollowUpId of 0 (they currently have no value
# in column H at all).
$defaultRows = @(3, 5, 6, 8, 10, 11, 12, 13, 15, 16, 17, 18, 20, 21, 22, 23, 24, 25, 27, 29)
foreach ($r in $defaultRows) {
    $ws.Cells.Item($r, 8).Value = 0
}

# Leave the existing links alone (H4 = 2000, H7 = 2001, H28 = 2002).

# Reflect the resulting selection / view state.
$ws.Range("H28").Select() | Out-Null
